# BulkImportDateParsing.xlsx - add parsed Date/Time columns, fix "sixth" -> "Sixth",
# rename the raw-string date column header to "DateTimeCol".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixups on existing columns -------------------------------------
$ws.Range("B1").Value = "DateTimeCol"
$ws.Range("A5").Value = "Sixth"

# --- New header row cells --------------------------------------------------
$ws.Range("C1").Value = "DateCol"
$ws.Range("D1").Value = "TimeCol"

# --- Column C: parsed date (whole-day serials), format mm-dd-yy (builtin 14) --
$ws.Range("C2").Value = 44188
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)
$ws.Range("C3").Value = 43853
$ws.Range("C4").Value = 43913
$ws.Range("C5").Value = 44005

# --- Column D: parsed time-of-day fraction, custom AM/PM format --------------
$ws.Range("D2").Value = 0.51428240740740738
$ws.Range("D2").NumberFormat = "[`$-F400]h:mm:ss\ AM/PM"
$ws.Range("D2").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$ws.Range("D3").Value = 0.51706018518518515
$ws.Range("D4").Value = 0.51636574074074071
$ws.Range("D5").Value = 0.52053240740740747

$ws.Application.CutCopyMode = $false

# --- Column widths (characters) ---------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.330729166666666
$ws.Columns.Item(2).ColumnWidth = 31.498697916666664
$ws.Columns.Item(4).ColumnWidth = 8.998697916666666

# --- Normalize row heights (drop explicit 13.5pt overrides) ----------------
$ws.Range("A2:A5").EntireRow.AutoFit()

# --- Selection matches the authored file ------------------------------------
$ws.Range("D2").Select()
